# Aula4 "Produtos Atualizados" - update currency quotations (Cotacao) and
# the derived "Preco Base Reais" / "Preco Final" columns after pulling the
# new rates (cf. commit: "Adicionando projeto de pegar cotacoes via yahoo").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Cotação) / Column E (Preço Base Reais) ------------------
$ws.Range("D2").Value = 5.4346
$ws.Range("E2").Value = 5434.545654

$ws.Range("D3").Value = 6.5363
$ws.Range("E3").Value = 29413.35

$ws.Range("D4").Value = 5.4346
$ws.Range("E4").Value = 4891.085654

$ws.Range("D5").Value = 5.4346
$ws.Range("E5").Value = 4342.2454

$ws.Range("D6").Value = 6.5363
$ws.Range("E6").Value = 19608.9

$ws.Range("D7").Value = 5.4346
$ws.Range("E7").Value = 2611.216608

$ws.Range("D8").Value = 308.59
$ws.Range("E8").Value = 6171.799999999999

# --- Column G (Preço Final) ---------------------------------------------
# These are stored as text in the workbook (shared strings), not numbers,
# so the cell's number format is forced to Text before the write and then
# restored to the default style afterwards (keeps the same style index
# the cells had originally).
$ws.Range("G2:G8").ClearContents()

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "7608.36"
$ws.Range("G2").Style = "Normal"

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "58826.70"
$ws.Range("G3").Style = "Normal"

$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "8314.85"
$ws.Range("G4").Style = "Normal"

$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "7381.82"
$ws.Range("G5").Style = "Normal"

$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "37256.91"
$ws.Range("G6").Style = "Normal"

$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "5222.43"
$ws.Range("G7").Style = "Normal"

$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "7097.57"
$ws.Range("G8").Style = "Normal"
